# edit.ps1 -- apply the "Context" section rewrite to linearityReportMarch2013Data.docx
#
# Summary of the change:
#   1. Paragraph "The Clio infrared camera..." - the middle of the paragraph
#      (up to "...sensitive from 1 ") is rewritten/condensed into a single run,
#      leaving the trailing "micro-m ... micro-m." runs untouched.
#   2. Paragraph "The context to this report is to..." is replaced wholesale
#      with new wording about the report's goal / timeline.
#   3. The two paragraphs that used to read "This data was originally
#      gathered..." and "Now, linearity is what we want to achieve..." swap
#      roles: the (rewritten) linearity paragraph now comes first, followed
#      by the (rewritten, new date) "originally gathered" paragraph.
#   4. In the trailing empty paragraph, the `_GoBack` bookmark now appears
#      before the page-break run instead of after it.

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find.Execute could not locate target text: $oldText"
    }
}

# --- 1. "The Clio infrared camera..." intro paragraph -----------------------
$old1 = "The Clio infrared camera is an instrument used in the Magellan Adaptive Optics system, located at the Las Campanas Observatory located in the Atacama Desert in Chile. Since the Earth" + [char]0x2019 + "s atmosphere is mostly opaque to infrared light, infrared astronomy should be conducted at high elevations, past as much of the atmosphere as possible. Also, adaptive optics are used to gather data and correct for the turbulence of the atmosphere. However, Clio is specifically sensitive from 1 "
$new1 = "The Clio infrared camera is an instrument used in the Magellan Adaptive Optics system, located at the Las Campanas Observatory located in the Atacama Desert in Chile. Clio specialized in infrared photometry, and adaptive optics are used to ensure that the turbulence of the atmosphere are corrected for and clearer pictures can be taken. Specifically, Clio is sensitive from 1 "
Replace-Text $old1 $new1

# --- 2. "The context to this report is to..." paragraph ---------------------
$old2 = "The context to this report is to analyze a set of images to calibrate the infrared camera Clio.  To do this, a program must be written to judge how exactly the images should be corrected for linearity. Now, this data has previously been corrected for linearity by Katie Morzinski. However, I chose to also attempt to correct this data, as I want to use what I have done for this specific set of data as a template for correcting data that hasn" + [char]0x2019 + "t been corrected yet."
$new2 = "Now, the goal of this report was to analyze a set of images to calibrate the infrared camera CLIO.  The code to calibrate this data set was written from August 2016 " + [char]0x2013 + " January 2017 and adapted from code that I have previously wrote to analyze another data set from the Clio camera."
Replace-Text $old2 $new2

# --- 3. Swap + rewrite the "originally gathered" / "linearity" paragraphs ---
# The paragraph that used to be "This data was originally gathered..." now
# holds the (new) linearity discussion, and the paragraph that used to hold
# the linearity discussion now holds the (new) "originally gathered" text.
$old3 = "This data was originally gathered on March 23, 2013, with the Clio camera in the MagAO system at the Las Campanas observatory in Chile."
$new3 = "In the end, linearity is what we want to achieve with the data set. The supposed relationship in the data between the integration time and counts readings should be linear, however, due to saturation from increased brightness, non-linearity appears within the data trends. Thus, the data is rendered useless. It is therefore the objective of this report to demonstrate that through curve fitting and other programming tools, we can correct this data and make sure that a larger portion of it turns out to be linear. This is known as linearity, and it is crucial to preserving the viability of the data."
Replace-Text $old3 $new3

$old4 = "Now, linearity is what we want to achieve with the data set. The supposed relationship between the integration time and counts readings should be linear, however, due to external factors, the higher the integration time, the more counts are dropped from the reading. Thus, the data is rendered useless. It is therefore the goal of this report to demonstrate that through curve fitting and other programming tools, we can correct this data and make sure that a larger portion of it turns out to be linear, making the data useful. This is known as linearity, and it is crucial to the correction of this data."
$new4 = "This data was originally gathered in November of 2014, with the Clio camera in the MagAO system at the Las Campanas observatory in Chile."
Replace-Text $old4 $new4

# --- 4. Move the `_GoBack` bookmark ahead of the page-break run ------------
$bm = $d.Bookmarks("_GoBack")
$paraStart = $bm.Range.Paragraphs(1).Range.Start
$bm.Delete()
$newBmRange = $d.Range($paraStart, $paraStart)
$d.Bookmarks.Add("_GoBack", $newBmRange)
